# Updated symbol list with latest coin ranking data.
# All Price/Volume/Coin/Link cells are stored as text in this sheet,
# so we force a Text number format while writing to preserve exact
# formatting (e.g. trailing zeros like "0.001250"), then restore the
# original "General" number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "D19", "E19", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "B48", "C48", "D48", "E48", "B49", "C49", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($cellRef in $targetCells) { $ws.Range($cellRef).NumberFormat = "@" }

$ws.Range("D2").Value = "307.81"
$ws.Range("E2").Value = "-1.34%"
$ws.Range("D3").Value = "40.87"
$ws.Range("E3").Value = "0.57%"
$ws.Range("D4").Value = "5.044"
$ws.Range("E4").Value = "-1.45%"
$ws.Range("D5").Value = "0.07638"
$ws.Range("E5").Value = "-2.73%"
$ws.Range("D6").Value = "4.261"
$ws.Range("E6").Value = "-1.82%"
$ws.Range("D7").Value = "1.611"
$ws.Range("E7").Value = "-3.71%"
$ws.Range("D8").Value = "0.9086"
$ws.Range("E8").Value = "-1.69%"
$ws.Range("D9").Value = "2.444"
$ws.Range("E9").Value = "-4.30%"
$ws.Range("D10").Value = "0.1011"
$ws.Range("E10").Value = "-6.01%"
$ws.Range("D11").Value = "0.1773"
$ws.Range("E11").Value = "-0.88%"
$ws.Range("D12").Value = "0.09167"
$ws.Range("E12").Value = "0.79%"
$ws.Range("D13").Value = "0.04232"
$ws.Range("E13").Value = "-4.90%"
$ws.Range("D14").Value = "0.1053"
$ws.Range("E14").Value = "-0.43%"
$ws.Range("D15").Value = "0.001250"
$ws.Range("E15").Value = "-1.14%"
$ws.Range("D16").Value = "0.005884"
$ws.Range("E16").Value = "-0.11%"
$ws.Range("D17").Value = "3.362"
$ws.Range("E17").Value = "0.23%"
$ws.Range("E18").Value = "-1.50%"
$ws.Range("D19").Value = "6.768"
$ws.Range("E19").Value = "-5.94%"
$ws.Range("E20").Value = "-1.77%"
$ws.Range("D21").Value = "0.2722"
$ws.Range("E21").Value = "2.52%"
$ws.Range("D22").Value = "0.04155"
$ws.Range("E22").Value = "-0.73%"
$ws.Range("D23").Value = "0.001217"
$ws.Range("E23").Value = "-2.46%"
$ws.Range("D24").Value = "0.004078"
$ws.Range("E24").Value = "-1.64%"
$ws.Range("D25").Value = "0.0001300"
$ws.Range("E25").Value = "5.44%"
$ws.Range("D26").Value = "0.0003006"
$ws.Range("E26").Value = "-0.01%"
$ws.Range("D38").Value = "0.02416"
$ws.Range("E38").Value = "-1.21%"
$ws.Range("D39").Value = "0.05181"
$ws.Range("E39").Value = "-1.96%"
$ws.Range("D40").Value = "0.007783"
$ws.Range("E40").Value = "-3.44%"
$ws.Range("D41").Value = "0.1309"
$ws.Range("E41").Value = "-3.44%"
$ws.Range("D42").Value = "0.007091"
$ws.Range("E42").Value = "8.01%"
$ws.Range("D43").Value = "0.001948"
$ws.Range("E43").Value = "-4.40%"
$ws.Range("D44").Value = "0.007480"
$ws.Range("E44").Value = "-9.61%"
$ws.Range("D45").Value = "0.3057"
$ws.Range("E45").Value = "-1.69%"
$ws.Range("D46").Value = "0.00006377"
$ws.Range("E46").Value = "-6.30%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.94%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "0.005290"
$ws.Range("E48").Value = "54.15%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "0.004398"
$ws.Range("E49").Value = "6.30%"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "-0.94%"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "-0.94%"

foreach ($cellRef in $targetCells) { $ws.Range($cellRef).NumberFormat = "General" }
